$d = $word.ActiveDocument

# The document currently ends with an empty trailing paragraph. Turn it
# into a paragraph with two runs: a bold "Git reflog ->" label followed
# by a normal (non-bold) explanatory sentence -- matching the style used
# by every other "Git <cmd> -> explanation" paragraph already present.

$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$r = $p.Range

$labelText = "Git reflog ->"
$tailText  = " Vê as versões adicionas até o momento, ou seja, os commits."

# 1) Write the plain (non-bold) explanatory text into the paragraph first.
$r.Text = $tailText

# 2) Insert the label text in front of it. It starts out with the same
#    (no) formatting as the text that follows, so Word keeps it as a
#    single run for now.
$insertPoint = $d.Range($r.Start, $r.Start)
$insertPoint.InsertBefore($labelText)

$pStart = $d.Paragraphs.Item($count).Range.Start
$labelRange = $d.Range($pStart, $pStart + $labelText.Length)

# 3) Borrow the bold run-formatting (w:b + w:bCs) from one of the other
#    bold "Git ... ->" labels already in the document (the paragraph
#    right before this one), then restore the label's own text
#    (FormattedText assignment pulls in the source run's text too, so
#    we overwrite it back to "Git reflog ->" right after).
$srcPara = $d.Paragraphs.Item($count - 1)
$srcLabelLen = "Git push ->".Length
$srcBoldRange = $d.Range($srcPara.Range.Start, $srcPara.Range.Start + $srcLabelLen)

$labelRange.FormattedText = $srcBoldRange.FormattedText
$labelRange = $d.Range($pStart, $pStart + $srcLabelLen)
$labelRange.Text = $labelText
